$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (and formatting) to the right
$ws.Columns.Item(1).Insert()

# Copy the header formatting from the (now shifted) neighboring header cell B1 onto the new A1 header
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header and values for the "tipo_renta" column
$ws.Cells.Item(1, 1).Value = "tipo_renta"
$ws.Cells.Item(2, 1).Value = "FIJA"
$ws.Cells.Item(3, 1).Value = "VARIABLE"
